$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Normalize the old "last row" (row 7: Battery / 12 V) to the regular
#     row styling used by the other data rows (copy format from row 2). ---
$ws.Range("A2:C2").Copy() | Out-Null
$ws.Range("A7:C7").PasteSpecial(-4122) | Out-Null   # xlPasteFormats

# --- Add the new row 8 (UGV Kit), reusing the styling pattern that the
#     previous last row (row 7) used to have: bordered/centered cells in
#     columns A & C, and an empty bordered cell (like B6) in column B. ---
$ws.Range("A2").Copy() | Out-Null
$ws.Range("A8").PasteSpecial(-4122) | Out-Null

$ws.Range("C2").Copy() | Out-Null
$ws.Range("C8").PasteSpecial(-4122) | Out-Null

$ws.Range("B6").Copy() | Out-Null
$ws.Range("B8").PasteSpecial(-4122) | Out-Null

$excel.CutCopyMode = 0

$ws.Range("A8").Value = "UGV Kit"
$ws.Range("C8").Value = 1

# --- Update the saved selection to match the new active cell. ---
$ws.Range("H16").Select() | Out-Null
